$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12 label (2025-08-30) as TEXT, not auto-converted to a date ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A12").ClearFormats()

# --- Bulk numeric update for B2:K12 (rows 2-11 revised, row 12 new) ---
$data = New-Object 'object[,]' 11,10
$data[0,0] = 0.43733570941851607
$data[0,1] = 0
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 0
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0
$data[0,9] = 0
$data[1,0] = 0.38609305706739067
$data[1,1] = 0
$data[1,2] = -0.03771250672027778
$data[1,3] = 0.002167778819963516
$data[1,4] = -0.0065804770510286825
$data[1,5] = 0.02013926719455408
$data[1,6] = -0.0021194881224011093
$data[1,7] = -0.03582240899056864
$data[1,8] = 0
$data[1,9] = 0.008685182518633205
$data[2,0] = 0.386933510793528
$data[2,1] = 0.014949331827655529
$data[2,2] = 0
$data[2,3] = 0.006244366588589525
$data[2,4] = -0.004104656478545085
$data[2,5] = 0
$data[2,6] = 0.002741752272146474
$data[2,7] = -0.004829222293115547
$data[2,8] = -0.011507576218694268
$data[2,9] = -0.0026535419718992648
$data[3,0] = 0.4986762201898488
$data[3,1] = 0.15483817319727505
$data[3,2] = -0.010247081911041006
$data[3,3] = -0.013129856503409105
$data[3,4] = 0.00870610800843658
$data[3,5] = -0.04035755547330965
$data[3,6] = 0.00495020557835481
$data[3,7] = -0.007822302199183923
$data[3,8] = 0
$data[3,9] = 0.014805018699198036
$data[4,0] = 0.40223077763667914
$data[4,1] = -0.010355563066151244
$data[4,2] = 0
$data[4,3] = -0.1790242177454549
$data[4,4] = 0.0010346888043497218
$data[4,5] = 0
$data[4,6] = 0.01253899025945475
$data[4,7] = 0.06778608528400032
$data[4,8] = 0
$data[4,9] = 0.011574573910631658
$data[5,0] = 0.37530726784054225
$data[5,1] = 0
$data[5,2] = 0.06148563210584754
$data[5,3] = -0.0023420756961862565
$data[5,4] = -0.1058971774990902
$data[5,5] = -0.002535490048060517
$data[5,6] = 0
$data[5,7] = -0.0052581111710945675
$data[5,8] = 0
$data[5,9] = 0.027623712512447107
$data[6,0] = 0.12830453245437112
$data[6,1] = -0.23734721367592448
$data[6,2] = 0
$data[6,3] = 0.053723759091085725
$data[6,4] = 0.013147778961419087
$data[6,5] = 0
$data[6,6] = -0.0009480529160716372
$data[6,7] = -0.07436824079466713
$data[6,8] = 0
$data[6,9] = -0.0012107660520126462
$data[7,0] = 0.028058698140802363
$data[7,1] = 0
$data[7,2] = -0.059111275493806306
$data[7,3] = 0.0010374724631778101
$data[7,4] = 0.03519889810463913
$data[7,5] = -0.054182998919725264
$data[7,6] = -0.0013003166946010012
$data[7,7] = -0.013416413356474657
$data[7,8] = 0
$data[7,9] = -0.008471200416778485
$data[8,0] = 0.3267710828856211
$data[8,1] = 0.33112973181037714
$data[8,2] = 0
$data[8,3] = -0.014297089691692876
$data[8,4] = -0.010616102648969903
$data[8,5] = 0
$data[8,6] = -0.0006223601373289377
$data[8,7] = 0.011571330186565566
$data[8,8] = -0.001005151568616337
$data[8,9] = -0.017447973205515943
$data[9,0] = 0.6187348168057625
$data[9,1] = 0
$data[9,2] = 0.2957468444718051
$data[9,3] = -0.014021195040515633
$data[9,4] = 0.06599290665306784
$data[9,5] = -0.0665814579285608
$data[9,6] = -0.002617195194888335
$data[9,7] = 0.066427941433685
$data[9,8] = 0
$data[9,9] = -0.05298411047445184
$data[10,0] = 0.3778246626708002
$data[10,1] = -0.2192026031591393
$data[10,2] = 0
$data[10,3] = -0.0036759534487288584
$data[10,4] = -0.002046347821404191
$data[10,5] = 0
$data[10,6] = -0.0004392631243039243
$data[10,7] = -0.015287234045076557
$data[10,8] = 0
$data[10,9] = -0.000258752536309359
$ws.Range("B2:K12").Value = $data

# --- Column width tweaks (C,D narrower; E,F,J adjusted) ---
$ws.Columns.Item(3).ColumnWidth = 13.25
$ws.Columns.Item(4).ColumnWidth = 13.25
$ws.Columns.Item(5).ColumnWidth = 14.25
$ws.Columns.Item(6).ColumnWidth = 14.25
$ws.Columns.Item(10).ColumnWidth = 14.25

